$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seguimiento a riesgos")

# Populate the new risk entry in row 20 (order matters so new shared
# strings are appended in the same sequence as the authored workbook)
$ws.Range("A20").Value = 10
$ws.Range("B20").Value = "Los integrantes del equipo de desarrollo no administran de manera correcta el sofware de gestion de proyecto"
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Formula = "=D20*C20"
$ws.Range("F20").Value = "utilizar y actualizar el software de gestion para administrar el proyecto a diario"
$ws.Range("G20").Value = "Cambiar de software para llevar la gestion del seguimiento del proyecto"
$ws.Range("I20").Value = "Presentado"

# Fill the "Última fecha de revisión" (H) column for rows 11-20 with the date text
$fecha = "17/05/2019"
for ($r = 11; $r -le 20; $r++) {
    $ws.Cells.Item($r, 8).Value = $fecha
}

# Row 20 height change (grows to fit the wrapped risk description)
$ws.Rows.Item(20).RowHeight = 29.25

# Update the active selection shown in the sheet view
$ws.Range("F13").Select()
